# Update "想去人数" (people interested) counts on the "展览" and "全部类型"
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 77
$ws1.Range("F4").Value = 1164
$ws1.Range("F5").Value = 1332
$ws1.Range("F6").Value = 291
$ws1.Range("F7").Value = 1018
$ws1.Range("F8").Value = 10475
$ws1.Range("F9").Value = 11
$ws1.Range("F11").Value = 279
$ws1.Range("F14").Value = 11962
$ws1.Range("F15").Value = 12368
$ws1.Range("F17").Value = 112

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 77
$ws4.Range("F5").Value = 1164
$ws4.Range("F6").Value = 1332
$ws4.Range("F7").Value = 291
$ws4.Range("F8").Value = 1018
$ws4.Range("F9").Value = 10475
$ws4.Range("F10").Value = 11
$ws4.Range("F12").Value = 279
$ws4.Range("F15").Value = 11962
$ws4.Range("F16").Value = 12368
$ws4.Range("F18").Value = 112
